# Generate Report for Handoff
# Update the handoff-run identifiers (old run "12f82e37..." -> new run
# "4986cffc...") and their associated timestamps / handoff file names across
# the Overview, zh-cn and de-de sheets, including the display text shown on
# the hyperlinks that point at those file names.

$wb = $excel.ActiveWorkbook

$oldId = "12f82e37-4a56-4ab3-a1d4-5d3d4e81705a"
$newId = "4986cffc-a556-4a06-ba95-3d7c041b717c"

$oldZhXlf = "$oldId.b257fad54a254673392e55cba252043fb79d0a50.zh-cn.xlf"
$newZhXlf = "$newId.7f3fb16251e652cae5fa6441d5a2b3c20c7cca3a.zh-cn.xlf"

$oldDeXlf = "$oldId.b257fad54a254673392e55cba252043fb79d0a50.de-de.xlf"
$newDeXlf = "$newId.7f3fb16251e652cae5fa6441d5a2b3c20c7cca3a.de-de.xlf"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/7a7f2ca76896cd28743061cfa0d18c431b15fcd6/e2e/$oldId.md"
$zhAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a3c594ca2b6a1e547f2e2ec282df022a43dc334/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/098f0d697196abda2f27e0612e89504999c98852/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# Original hyperlink font colour (#6495ED, "cornflower blue") expressed as a
# VBA-style BGR long so re-adding a hyperlink keeps the workbook's existing
# look instead of falling back to Excel's default hyperlink blue.
$hyperlinkColor = 15570276

function Update-Hyperlink($ws, $cellRef, $address, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText)
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# --- Sheet: Overview ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "$newId.md"
$ws1.Range("D2").Value = "2016-03-23 11:07:41"

# Hyperlinks.Delete() clears every hyperlink on the sheet, so gather all of
# them, wipe the sheet's hyperlinks once, then re-add each one with its
# (possibly updated) display text.
$ws1.Cells.Hyperlinks.Delete()
Update-Hyperlink $ws1 "A2" $mdAddress "$newId.md"

# --- Sheet: zh-cn ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "$newId.md"
$ws2.Range("D2").Value = $newZhXlf
$ws2.Range("E2").Value = "2016-03-23 11:07:34"

$ws2.Cells.Hyperlinks.Delete()
Update-Hyperlink $ws2 "A2" $mdAddress "$newId.md"
Update-Hyperlink $ws2 "D2" $zhAddress $newZhXlf

# --- Sheet: de-de ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "$newId.md"
$ws3.Range("D2").Value = $newDeXlf
$ws3.Range("E2").Value = "2016-03-23 11:07:41"

$ws3.Cells.Hyperlinks.Delete()
Update-Hyperlink $ws3 "A2" $mdAddress "$newId.md"
Update-Hyperlink $ws3 "D2" $deAddress $newDeXlf
